$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has rows for 2008..2020 (row2..row14).
# Target: rows for 2010..2021 (row2..row13) -> drop the 2008/2009 rows
# (shifting everything up by two), then append a brand-new 2021 row
# with its data at the end.

# Remove the 2008 row; this shifts 2009..2020 up into rows 2..13
$ws.Rows.Item(2).Delete()
# Remove the (now shifted) 2009 row; this shifts 2010..2020 up into rows 2..12
$ws.Rows.Item(2).Delete()

# Append the new 2021 row (row 13) with its reported values.
# Row 13 doesn't exist yet after the deletes above, so clone the label
# cell's formatting (bold, centered, bordered) from the row above it.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 28.4
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = 221884.6
$ws.Range("F13").Value = 214036
$ws.Range("G13").Value = 304072.8
$ws.Range("H13").Value = 94560.89999999999
$ws.Range("I13").Value = 17851.8
$ws.Range("J13").Value = ""
$ws.Range("K13").Value = 116496.1
$ws.Range("L13").Value = 211197.7
$ws.Range("M13").Value = 113657.8
$ws.Range("N13").Value = ""
$ws.Range("O13").Value = ""
$ws.Range("P13").Value = ""
$ws.Range("Q13").Value = ""
$ws.Range("R13").Value = 5740.8
$ws.Range("S13").Value = ""
$ws.Range("T13").Value = 327693.8
$ws.Range("U13").Value = ""
$ws.Range("V13").Value = 11248.4
